# Apply the workbook data refresh: update the "F" column (售罄/剩余 count, etc.)
# numeric values that changed between the two generated snapshots.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (展览信息)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 38041
$ws1.Range("F11").Value = 768
$ws1.Range("F15").Value = 48
$ws1.Range("F22").Value = 895
$ws1.Range("F23").Value = 2619
$ws1.Range("F24").Value = 1108
$ws1.Range("F26").Value = 133
$ws1.Range("F29").Value = 861
$ws1.Range("F31").Value = 1202

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F9").Value = 146

# Sheet "全部类型" (combined listing, mirrors rows from the sheets above)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 38041
$ws4.Range("F18").Value = 768
$ws4.Range("F24").Value = 146
$ws4.Range("F26").Value = 48
$ws4.Range("F34").Value = 895
$ws4.Range("F35").Value = 2619
$ws4.Range("F36").Value = 1108
$ws4.Range("F38").Value = 133
$ws4.Range("F42").Value = 861
$ws4.Range("F44").Value = 1202
